$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "75.945.89"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "2.906.89"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "200.32"
$ws.Range("E5").Value = "  +6.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "596.24"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.549"
$ws.Range("E8").Value = "  -0.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.198"
$ws.Range("E9").Value = "  +2.99%  "
$ws.Range("D10").Value = "2.904.93"
$ws.Range("E10").Value = "  +2.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.429"
$ws.Range("E11").Value = "  +15.64%  "
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").Value = "3.440.21"
$ws.Range("E14").Value = "  +2.98%  "
$ws.Range("D15").Value = "75.839.84"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.62"
$ws.Range("E16").Value = "  +2.72%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000191"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").Value = "2.899.76"
$ws.Range("E18").Value = "  +2.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.92"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.83"
$ws.Range("E20").Value = "  +4.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "371.90"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.31"
$ws.Range("E22").Value = "  +2.73%  "
$ws.Range("E23").Value = "  +4.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.99"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "3.043.97"
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.18"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.64"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000108"
$ws.Range("E29").Value = "  +4.19%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "501.36"
$ws.Range("E32").Value = "  -2.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.69"
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.81"
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.96"
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.13"
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.65"
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.113"
$ws.Range("E39").Value = "  -5.29%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.101"
$ws.Range("E41").Value = "  +18.23%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "180.18"
$ws.Range("E42").Value = "  -2.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.348"
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.96"
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.08"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.33"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.569"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.71"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.28"
$ws.Range("E51").Value = "  +6.12%  "
